$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force all target cells to Text format so numeric-looking strings
# (e.g. "45.456.57", "15.50", "1.00") are preserved verbatim, matching
# the source data which stores these as literal text, not numbers.
$ws.Range('D2').Value = "'45.456.57"
$ws.Range('E2').Value = "'  -0.93%  "
$ws.Range('D3').Value = "'2.375.99"
$ws.Range('E3').Value = "'  -1.96%  "
$ws.Range('E4').Value = "'  -0.04%  "
$ws.Range('D5').Value = "'318.74"
$ws.Range('E5').Value = "'  -0.38%  "
$ws.Range('D6').Value = "'108.92"
$ws.Range('E6').Value = "'  -6.84%  "
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E8').Value = "'  +0.04%  "
$ws.Range('D9').Value = "'0.616"
$ws.Range('E9').Value = "'  -2.84%  "
$ws.Range('D10').Value = "'41.01"
$ws.Range('E10').Value = "'  -5.56%  "
$ws.Range('E11').Value = "'  -2.23%  "
$ws.Range('E12').Value = "'  -2.69%  "
$ws.Range('E13').Value = "'  +0.21%  "
$ws.Range('E14').Value = "'  -5.04%  "
$ws.Range('D15').Value = "'2.738.09"
$ws.Range('E15').Value = "'  -1.94%  "
$ws.Range('D16').Value = "'15.50"
$ws.Range('E16').Value = "'  -3.26%  "
$ws.Range('D17').Value = "'2.372.27"
$ws.Range('E17').Value = "'  -2.40%  "
$ws.Range('D18').Value = "'45.422.74"
$ws.Range('E18').Value = "'  -0.82%  "
$ws.Range('D19').Value = "'15.27"
$ws.Range('E19').Value = "'  +13.91%  "
$ws.Range('D20').Value = "'7.33"
$ws.Range('E20').Value = "'  -4.10%  "
$ws.Range('E21').Value = "'  -2.52%  "
$ws.Range('D22').Value = "'3.63"
$ws.Range('E22').Value = "'  +1.64%  "
$ws.Range('D23').Value = "'73.37"
$ws.Range('E23').Value = "'  -2.50%  "
$ws.Range('D24').Value = "'264.86"
$ws.Range('E24').Value = "'  -1.75%  "
$ws.Range('E25').Value = "'  -2.21%  "
$ws.Range('E26').Value = "'  +0.18%  "
$ws.Range('E27').Value = "'  -1.20%  "
$ws.Range('D28').Value = "'7.50"
$ws.Range('E28').Value = "'  -2.02%  "
$ws.Range('E29').Value = "'  -1.87%  "
$ws.Range('D30').Value = "'22.48"
$ws.Range('E30').Value = "'  -2.76%  "
$ws.Range('D31').Value = "'0.0958"
$ws.Range('E31').Value = "'  -0.73%  "
$ws.Range('D32').Value = "'37.29"
$ws.Range('E32').Value = "'  -7.66%  "
$ws.Range('D33').Value = "'169.16"
$ws.Range('E33').Value = "'  -2.80%  "
$ws.Range('E34').Value = "'  -4.19%  "
$ws.Range('B35').Value = "'LidoDAOToken"
$ws.Range('C35').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('D35').Value = "'3.32"
$ws.Range('E35').Value = "'  +5.19%  "
$ws.Range('B36').Value = "'Stellar"
$ws.Range('C36').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('D36').Value = "'0.132"
$ws.Range('E36').Value = "'  -0.43%  "
$ws.Range('E37').Value = "'  -2.73%  "
$ws.Range('D38').Value = "'4.75"
$ws.Range('E38').Value = "'  -5.16%  "
$ws.Range('D39').Value = "'1.98"
$ws.Range('E39').Value = "'  +8.05%  "
$ws.Range('D40').Value = "'4.02"
$ws.Range('E40').Value = "'  -4.44%  "
$ws.Range('E41').Value = "'  -3.45%  "
$ws.Range('D42').Value = "'98.20"
$ws.Range('E42').Value = "'  -4.21%  "
$ws.Range('D43').Value = "'70.75"
$ws.Range('E43').Value = "'  -2.43%  "
$ws.Range('B44').Value = "'Maker"
$ws.Range('C44').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D44').Value = "'1.881.79"
$ws.Range('E44').Value = "'  +13.10%  "
$ws.Range('B45').Value = "'Celestia"
$ws.Range('C45').Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range('D45').Value = "'13.04"
$ws.Range('E45').Value = "'  -5.00%  "
$ws.Range('D46').Value = "'0.230"
$ws.Range('E46').Value = "'  -4.52%  "
$ws.Range('D47').Value = "'6.06"
$ws.Range('E47').Value = "'  +3.13%  "
$ws.Range('B48').Value = "'ordi"
$ws.Range('C48').Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range('D48').Value = "'86.58"
$ws.Range('E48').Value = "'  +7.20%  "
$ws.Range('B49').Value = "'FirstDigitalUSD"
$ws.Range('C49').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('D49').Value = "'1.00"
$ws.Range('E49').Value = "'  +0.06%  "
$ws.Range('D50').Value = "'9.44"
$ws.Range('E50').Value = "'  -0.41%  "
$ws.Range('D51').Value = "'112.76"
$ws.Range('E51').Value = "'  -4.32%  "
